$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells as Text so numeric-looking strings (e.g. "1.001")
# are preserved verbatim instead of being parsed into floating point numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.918.52'
$ws.Range("E2").Value = '  -1.08%  '

$ws.Range("D3").Value = '1.951.44'
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '242.40'
$ws.Range("E5").Value = '  -2.60%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").Value = '0.4875'
$ws.Range("E7").Value = '  -0.36%  '

$ws.Range("D8").Value = '0.2941'
$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("D9").Value = '0.06927'
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  +1.27%  '

$ws.Range("D11").Value = '106.93'
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").Value = '1.971.15'
$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").Value = '0.07757'
$ws.Range("E13").Value = '  -0.45%  '

$ws.Range("D14").Value = '5.341'
$ws.Range("E14").Value = '  -2.27%  '

$ws.Range("D15").Value = '0.6952'
$ws.Range("E15").Value = '  -2.10%  '

$ws.Range("D16").Value = '278.44'
$ws.Range("E16").Value = '  -2.85%  '

$ws.Range("D17").Value = '30.932.49'
$ws.Range("E17").Value = '  -1.07%  '

$ws.Range("D18").Value = '0.000007738'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").Value = '13.19'
$ws.Range("E19").Value = '  -1.27%  '

$ws.Range("D20").Value = '2.212.90'
$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '5.472'
$ws.Range("E22").Value = '  -3.07%  '

$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").Value = '6.486'
$ws.Range("E24").Value = '  -2.57%  '

$ws.Range("D25").Value = '9.735'
$ws.Range("E25").Value = '  -3.20%  '

$ws.Range("D26").Value = '167.96'
$ws.Range("E26").Value = '  -1.32%  '

$ws.Range("D27").Value = '19.68'
$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").Value = '2.167'
$ws.Range("E28").Value = '  -1.55%  '

$ws.Range("E29").Value = '  -2.88%  '

$ws.Range("D30").Value = '1.398'
$ws.Range("E30").Value = '  -3.29%  '

$ws.Range("D31").Value = '4.574'
$ws.Range("E31").Value = '  -5.97%  '

$ws.Range("D32").Value = '1.555'
$ws.Range("E32").Value = '  -2.88%  '

$ws.Range("D33").Value = '4.378'
$ws.Range("E33").Value = '  -3.36%  '

$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("D35").Value = '0.7516'
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("D36").Value = '1.162'
$ws.Range("E36").Value = '  -1.09%  '

$ws.Range("D37").Value = '2.732'
$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("E39").Value = '  -2.09%  '

$ws.Range("D40").Value = '6.494'
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").Value = '77.25'
$ws.Range("E41").Value = '  +4.02%  '

$ws.Range("D42").Value = '2.096'
$ws.Range("E42").Value = '  -1.91%  '

$ws.Range("D43").Value = '0.8982'
$ws.Range("E43").Value = '  +0.92%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4423'
$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '108.32'
$ws.Range("E45").Value = '  -1.82%  '

$ws.Range("D46").Value = '0.9996'
$ws.Range("E46").Value = '  -0.36%  '

$ws.Range("D47").Value = '7.723'
$ws.Range("E47").Value = '  +2.54%  '

$ws.Range("D48").Value = '995.17'
$ws.Range("E48").Value = '  +0.63%  '

$ws.Range("E49").Value = '  -2.21%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '35.75'
$ws.Range("E50").Value = '  -1.12%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.167'
$ws.Range("E51").Value = '  -2.78%  '
